# Auto-generated Excel COM-interop edit script implementing the commit diff:
#  - Insert a new sheet '삼양식품' immediately before '엔씨소프트' (full new dataset, rows 2-101)
#  - On '카카오', 'NAVER', '농심' and '엔씨소프트': fix the B105 value (was a 0 placeholder)
#    and append 5 new daily rows (106-110) of date/remn_amt data.
#
# NOTE: inserting a worksheet shifts the index-based identity of any worksheet
# reference captured at/after the insertion point, so the '엔씨소프트' worksheet
# is re-fetched by name (stable) immediately after it is pushed back a slot.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- 카카오 (sheet 1): fix row 105 and append rows 106-110 ---
$ws1.Cells.Item(105, 2).Value = 884360
$ws1.Cells.Item(106, 1).Value = 45971
$ws1.Cells.Item(106, 2).Value = 884327
$ws1.Cells.Item(107, 1).Value = 45972
$ws1.Cells.Item(107, 2).Value = 919516
$ws1.Cells.Item(108, 1).Value = 45973
$ws1.Cells.Item(108, 2).Value = 872304
$ws1.Cells.Item(109, 1).Value = 45974
$ws1.Cells.Item(109, 2).Value = 895098
$ws1.Cells.Item(110, 1).Value = 45975
$ws1.Cells.Item(110, 2).Value = 835105
$ws1.Range("A105").Copy()
$ws1.Range("A106:A110").PasteSpecial($xlPasteFormats)

# --- NAVER (sheet 2): fix row 105 and append rows 106-110 ---
$ws2.Cells.Item(105, 2).Value = 1340259
$ws2.Cells.Item(106, 1).Value = 45971
$ws2.Cells.Item(106, 2).Value = 1343051
$ws2.Cells.Item(107, 1).Value = 45972
$ws2.Cells.Item(107, 2).Value = 1377485
$ws2.Cells.Item(108, 1).Value = 45973
$ws2.Cells.Item(108, 2).Value = 1376852
$ws2.Cells.Item(109, 1).Value = 45974
$ws2.Cells.Item(109, 2).Value = 1356725
$ws2.Cells.Item(110, 1).Value = 45975
$ws2.Cells.Item(110, 2).Value = 1248374
$ws2.Range("A105").Copy()
$ws2.Range("A106:A110").PasteSpecial($xlPasteFormats)

# --- 농심 (sheet 3): fix row 105 and append rows 106-110 ---
$ws3.Cells.Item(105, 2).Value = 121095
$ws3.Cells.Item(106, 1).Value = 45971
$ws3.Cells.Item(106, 2).Value = 122966
$ws3.Cells.Item(107, 1).Value = 45972
$ws3.Cells.Item(107, 2).Value = 121414
$ws3.Cells.Item(108, 1).Value = 45973
$ws3.Cells.Item(108, 2).Value = 120949
$ws3.Cells.Item(109, 1).Value = 45974
$ws3.Cells.Item(109, 2).Value = 122408
$ws3.Cells.Item(110, 1).Value = 45975
$ws3.Cells.Item(110, 2).Value = 112458
$ws3.Range("A105").Copy()
$ws3.Range("A106:A110").PasteSpecial($xlPasteFormats)

# --- Insert new sheet '삼양식품' immediately before '엔씨소프트' ---
$wsNcsoftOrig = $wb.Worksheets.Item(4)
$wsSamyang = $wb.Worksheets.Add($wsNcsoftOrig)
$wsSamyang.Name = "삼양식품"

# Re-fetch '엔씨소프트' by name now that insertion shifted it to index 5
$wsNcsoft = $wb.Worksheets.Item("엔씨소프트")

# Header row
$wsSamyang.Cells.Item(1, 1).Value = "date"
$wsSamyang.Cells.Item(1, 2).Value = "remn_amt"
$wsNcsoft.Range("A1:B1").Copy()
$wsSamyang.Range("A1:B1").PasteSpecial($xlPasteFormats)

# Data rows 2-101
$wsSamyang.Cells.Item(2, 1).Value = 45828
$wsSamyang.Cells.Item(2, 2).Value = 399273
$wsSamyang.Cells.Item(3, 1).Value = 45831
$wsSamyang.Cells.Item(3, 2).Value = 408211
$wsSamyang.Cells.Item(4, 1).Value = 45832
$wsSamyang.Cells.Item(4, 2).Value = 392950
$wsSamyang.Cells.Item(5, 1).Value = 45833
$wsSamyang.Cells.Item(5, 2).Value = 381677
$wsSamyang.Cells.Item(6, 1).Value = 45834
$wsSamyang.Cells.Item(6, 2).Value = 385015
$wsSamyang.Cells.Item(7, 1).Value = 45835
$wsSamyang.Cells.Item(7, 2).Value = 383111
$wsSamyang.Cells.Item(8, 1).Value = 45838
$wsSamyang.Cells.Item(8, 2).Value = 362383
$wsSamyang.Cells.Item(9, 1).Value = 45839
$wsSamyang.Cells.Item(9, 2).Value = 358188
$wsSamyang.Cells.Item(10, 1).Value = 45840
$wsSamyang.Cells.Item(10, 2).Value = 348406
$wsSamyang.Cells.Item(11, 1).Value = 45841
$wsSamyang.Cells.Item(11, 2).Value = 340890
$wsSamyang.Cells.Item(12, 1).Value = 45842
$wsSamyang.Cells.Item(12, 2).Value = 334454
$wsSamyang.Cells.Item(13, 1).Value = 45845
$wsSamyang.Cells.Item(13, 2).Value = 354397
$wsSamyang.Cells.Item(14, 1).Value = 45846
$wsSamyang.Cells.Item(14, 2).Value = 383996
$wsSamyang.Cells.Item(15, 1).Value = 45847
$wsSamyang.Cells.Item(15, 2).Value = 419990
$wsSamyang.Cells.Item(16, 1).Value = 45848
$wsSamyang.Cells.Item(16, 2).Value = 392840
$wsSamyang.Cells.Item(17, 1).Value = 45849
$wsSamyang.Cells.Item(17, 2).Value = 398197
$wsSamyang.Cells.Item(18, 1).Value = 45852
$wsSamyang.Cells.Item(18, 2).Value = 381192
$wsSamyang.Cells.Item(19, 1).Value = 45853
$wsSamyang.Cells.Item(19, 2).Value = 382606
$wsSamyang.Cells.Item(20, 1).Value = 45854
$wsSamyang.Cells.Item(20, 2).Value = 378815
$wsSamyang.Cells.Item(21, 1).Value = 45855
$wsSamyang.Cells.Item(21, 2).Value = 343739
$wsSamyang.Cells.Item(22, 1).Value = 45856
$wsSamyang.Cells.Item(22, 2).Value = 312748
$wsSamyang.Cells.Item(23, 1).Value = 45859
$wsSamyang.Cells.Item(23, 2).Value = 319175
$wsSamyang.Cells.Item(24, 1).Value = 45860
$wsSamyang.Cells.Item(24, 2).Value = 311450
$wsSamyang.Cells.Item(25, 1).Value = 45861
$wsSamyang.Cells.Item(25, 2).Value = 324789
$wsSamyang.Cells.Item(26, 1).Value = 45862
$wsSamyang.Cells.Item(26, 2).Value = 323229
$wsSamyang.Cells.Item(27, 1).Value = 45863
$wsSamyang.Cells.Item(27, 2).Value = 327777
$wsSamyang.Cells.Item(28, 1).Value = 45866
$wsSamyang.Cells.Item(28, 2).Value = 321580
$wsSamyang.Cells.Item(29, 1).Value = 45867
$wsSamyang.Cells.Item(29, 2).Value = 325458
$wsSamyang.Cells.Item(30, 1).Value = 45868
$wsSamyang.Cells.Item(30, 2).Value = 312229
$wsSamyang.Cells.Item(31, 1).Value = 45869
$wsSamyang.Cells.Item(31, 2).Value = 338290
$wsSamyang.Cells.Item(32, 1).Value = 45870
$wsSamyang.Cells.Item(32, 2).Value = 331328
$wsSamyang.Cells.Item(33, 1).Value = 45873
$wsSamyang.Cells.Item(33, 2).Value = 328102
$wsSamyang.Cells.Item(34, 1).Value = 45874
$wsSamyang.Cells.Item(34, 2).Value = 329759
$wsSamyang.Cells.Item(35, 1).Value = 45875
$wsSamyang.Cells.Item(35, 2).Value = 332494
$wsSamyang.Cells.Item(36, 1).Value = 45876
$wsSamyang.Cells.Item(36, 2).Value = 335131
$wsSamyang.Cells.Item(37, 1).Value = 45877
$wsSamyang.Cells.Item(37, 2).Value = 339676
$wsSamyang.Cells.Item(38, 1).Value = 45880
$wsSamyang.Cells.Item(38, 2).Value = 355615
$wsSamyang.Cells.Item(39, 1).Value = 45881
$wsSamyang.Cells.Item(39, 2).Value = 343819
$wsSamyang.Cells.Item(40, 1).Value = 45882
$wsSamyang.Cells.Item(40, 2).Value = 332178
$wsSamyang.Cells.Item(41, 1).Value = 45883
$wsSamyang.Cells.Item(41, 2).Value = 319782
$wsSamyang.Cells.Item(42, 1).Value = 45887
$wsSamyang.Cells.Item(42, 2).Value = 318901
$wsSamyang.Cells.Item(43, 1).Value = 45888
$wsSamyang.Cells.Item(43, 2).Value = 322134
$wsSamyang.Cells.Item(44, 1).Value = 45889
$wsSamyang.Cells.Item(44, 2).Value = 302183
$wsSamyang.Cells.Item(45, 1).Value = 45890
$wsSamyang.Cells.Item(45, 2).Value = 301770
$wsSamyang.Cells.Item(46, 1).Value = 45891
$wsSamyang.Cells.Item(46, 2).Value = 318928
$wsSamyang.Cells.Item(47, 1).Value = 45894
$wsSamyang.Cells.Item(47, 2).Value = 327600
$wsSamyang.Cells.Item(48, 1).Value = 45895
$wsSamyang.Cells.Item(48, 2).Value = 362643
$wsSamyang.Cells.Item(49, 1).Value = 45896
$wsSamyang.Cells.Item(49, 2).Value = 368604
$wsSamyang.Cells.Item(50, 1).Value = 45897
$wsSamyang.Cells.Item(50, 2).Value = 380595
$wsSamyang.Cells.Item(51, 1).Value = 45898
$wsSamyang.Cells.Item(51, 2).Value = 371477
$wsSamyang.Cells.Item(52, 1).Value = 45901
$wsSamyang.Cells.Item(52, 2).Value = 378666
$wsSamyang.Cells.Item(53, 1).Value = 45902
$wsSamyang.Cells.Item(53, 2).Value = 381118
$wsSamyang.Cells.Item(54, 1).Value = 45903
$wsSamyang.Cells.Item(54, 2).Value = 386688
$wsSamyang.Cells.Item(55, 1).Value = 45904
$wsSamyang.Cells.Item(55, 2).Value = 378468
$wsSamyang.Cells.Item(56, 1).Value = 45905
$wsSamyang.Cells.Item(56, 2).Value = 371630
$wsSamyang.Cells.Item(57, 1).Value = 45908
$wsSamyang.Cells.Item(57, 2).Value = 372378
$wsSamyang.Cells.Item(58, 1).Value = 45909
$wsSamyang.Cells.Item(58, 2).Value = 340353
$wsSamyang.Cells.Item(59, 1).Value = 45910
$wsSamyang.Cells.Item(59, 2).Value = 345246
$wsSamyang.Cells.Item(60, 1).Value = 45911
$wsSamyang.Cells.Item(60, 2).Value = 371847
$wsSamyang.Cells.Item(61, 1).Value = 45912
$wsSamyang.Cells.Item(61, 2).Value = 363922
$wsSamyang.Cells.Item(62, 1).Value = 45915
$wsSamyang.Cells.Item(62, 2).Value = 358233
$wsSamyang.Cells.Item(63, 1).Value = 45916
$wsSamyang.Cells.Item(63, 2).Value = 348172
$wsSamyang.Cells.Item(64, 1).Value = 45917
$wsSamyang.Cells.Item(64, 2).Value = 348733
$wsSamyang.Cells.Item(65, 1).Value = 45918
$wsSamyang.Cells.Item(65, 2).Value = 361364
$wsSamyang.Cells.Item(66, 1).Value = 45919
$wsSamyang.Cells.Item(66, 2).Value = 376191
$wsSamyang.Cells.Item(67, 1).Value = 45922
$wsSamyang.Cells.Item(67, 2).Value = 373961
$wsSamyang.Cells.Item(68, 1).Value = 45923
$wsSamyang.Cells.Item(68, 2).Value = 365016
$wsSamyang.Cells.Item(69, 1).Value = 45924
$wsSamyang.Cells.Item(69, 2).Value = 376606
$wsSamyang.Cells.Item(70, 1).Value = 45925
$wsSamyang.Cells.Item(70, 2).Value = 357990
$wsSamyang.Cells.Item(71, 1).Value = 45926
$wsSamyang.Cells.Item(71, 2).Value = 336907
$wsSamyang.Cells.Item(72, 1).Value = 45929
$wsSamyang.Cells.Item(72, 2).Value = 345605
$wsSamyang.Cells.Item(73, 1).Value = 45930
$wsSamyang.Cells.Item(73, 2).Value = 352982
$wsSamyang.Cells.Item(74, 1).Value = 45931
$wsSamyang.Cells.Item(74, 2).Value = 363509
$wsSamyang.Cells.Item(75, 1).Value = 45932
$wsSamyang.Cells.Item(75, 2).Value = 355896
$wsSamyang.Cells.Item(76, 1).Value = 45940
$wsSamyang.Cells.Item(76, 2).Value = 342103
$wsSamyang.Cells.Item(77, 1).Value = 45943
$wsSamyang.Cells.Item(77, 2).Value = 323791
$wsSamyang.Cells.Item(78, 1).Value = 45944
$wsSamyang.Cells.Item(78, 2).Value = 315268
$wsSamyang.Cells.Item(79, 1).Value = 45945
$wsSamyang.Cells.Item(79, 2).Value = 320009
$wsSamyang.Cells.Item(80, 1).Value = 45946
$wsSamyang.Cells.Item(80, 2).Value = 309390
$wsSamyang.Cells.Item(81, 1).Value = 45947
$wsSamyang.Cells.Item(81, 2).Value = 302922
$wsSamyang.Cells.Item(82, 1).Value = 45950
$wsSamyang.Cells.Item(82, 2).Value = 353564
$wsSamyang.Cells.Item(83, 1).Value = 45951
$wsSamyang.Cells.Item(83, 2).Value = 311931
$wsSamyang.Cells.Item(84, 1).Value = 45952
$wsSamyang.Cells.Item(84, 2).Value = 297880
$wsSamyang.Cells.Item(85, 1).Value = 45953
$wsSamyang.Cells.Item(85, 2).Value = 295438
$wsSamyang.Cells.Item(86, 1).Value = 45954
$wsSamyang.Cells.Item(86, 2).Value = 296621
$wsSamyang.Cells.Item(87, 1).Value = 45957
$wsSamyang.Cells.Item(87, 2).Value = 294470
$wsSamyang.Cells.Item(88, 1).Value = 45958
$wsSamyang.Cells.Item(88, 2).Value = 287977
$wsSamyang.Cells.Item(89, 1).Value = 45959
$wsSamyang.Cells.Item(89, 2).Value = 272080
$wsSamyang.Cells.Item(90, 1).Value = 45960
$wsSamyang.Cells.Item(90, 2).Value = 289544
$wsSamyang.Cells.Item(91, 1).Value = 45961
$wsSamyang.Cells.Item(91, 2).Value = 286219
$wsSamyang.Cells.Item(92, 1).Value = 45964
$wsSamyang.Cells.Item(92, 2).Value = 268234
$wsSamyang.Cells.Item(93, 1).Value = 45965
$wsSamyang.Cells.Item(93, 2).Value = 279148
$wsSamyang.Cells.Item(94, 1).Value = 45966
$wsSamyang.Cells.Item(94, 2).Value = 279945
$wsSamyang.Cells.Item(95, 1).Value = 45967
$wsSamyang.Cells.Item(95, 2).Value = 274300
$wsSamyang.Cells.Item(96, 1).Value = 45968
$wsSamyang.Cells.Item(96, 2).Value = 273894
$wsSamyang.Cells.Item(97, 1).Value = 45971
$wsSamyang.Cells.Item(97, 2).Value = 274586
$wsSamyang.Cells.Item(98, 1).Value = 45972
$wsSamyang.Cells.Item(98, 2).Value = 269812
$wsSamyang.Cells.Item(99, 1).Value = 45973
$wsSamyang.Cells.Item(99, 2).Value = 270550
$wsSamyang.Cells.Item(100, 1).Value = 45974
$wsSamyang.Cells.Item(100, 2).Value = 299492
$wsSamyang.Cells.Item(101, 1).Value = 45975
$wsSamyang.Cells.Item(101, 2).Value = 275854
$wsNcsoft.Range("A2").Copy()
$wsSamyang.Range("A2:A101").PasteSpecial($xlPasteFormats)

# --- 엔씨소프트 (existing sheet, now position 5): fix row 105 and append rows 106-110 ---
$wsNcsoft.Cells.Item(105, 2).Value = 164764
$wsNcsoft.Cells.Item(106, 1).Value = 45971
$wsNcsoft.Cells.Item(106, 2).Value = 173526
$wsNcsoft.Cells.Item(107, 1).Value = 45972
$wsNcsoft.Cells.Item(107, 2).Value = 175556
$wsNcsoft.Cells.Item(108, 1).Value = 45973
$wsNcsoft.Cells.Item(108, 2).Value = 190467
$wsNcsoft.Cells.Item(109, 1).Value = 45974
$wsNcsoft.Cells.Item(109, 2).Value = 179884
$wsNcsoft.Cells.Item(110, 1).Value = 45975
$wsNcsoft.Cells.Item(110, 2).Value = 169151
$wsNcsoft.Range("A105").Copy()
$wsNcsoft.Range("A106:A110").PasteSpecial($xlPasteFormats)

